$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Clone row 6's formatting down into the four new rows (7-10) ---
$ws.Range("A6:G6").Copy()
$ws.Range("A7:G10").PasteSpecial(-4122)   # xlPasteFormats

# --- Populate the new contact rows ---
# Host + the two guests were entered first (names, then e-mails), the
# co-host row was added afterwards - this mirrors the order the new
# shared strings were actually appended in.
$ws.Range("A7").Value = "sushanthost"
$ws.Range("A9").Value = "sushantguest1"
$ws.Range("A10").Value = "sushantguest2"
$ws.Range("C7").Value = "sushanthost@test.com"
$ws.Range("C9").Value = "sushantguest1@test.com"
$ws.Range("C10").Value = "sushantguest2@test.com"
$ws.Range("A8").Value = "sushantcohost"
$ws.Range("C8").Value = "sushantcohost@test.com"

$ws.Range("B7").Value = "abc"
$ws.Range("B8").Value = "abc"
$ws.Range("B9").Value = "abc"
$ws.Range("B10").Value = "abc"

$ws.Range("D7").Value = "972 BRAHMS CT"
$ws.Range("D8").Value = "972 BRAHMS CT"
$ws.Range("D9").Value = "972 BRAHMS CT"
$ws.Range("D10").Value = "972 BRAHMS CT"

$ws.Range("E7").Value = "TROY"
$ws.Range("E8").Value = "TROY"
$ws.Range("E9").Value = "TROY"
$ws.Range("E10").Value = "TROY"

$ws.Range("F7").Value = "Michigan"
$ws.Range("F8").Value = "Michigan"
$ws.Range("F9").Value = "Michigan"
$ws.Range("F10").Value = "Michigan"

$ws.Range("G7").Value = 48085
$ws.Range("G8").Value = 48085
$ws.Range("G9").Value = 48085
$ws.Range("G10").Value = 48085

# Match row 6's explicit row height
$ws.Range("A7:G10").RowHeight = 15

# --- Hyperlink the e-mail addresses, in the order they were added ---
$ws.Hyperlinks.Add($ws.Range("C7"), "mailto:sushanthost@test.com")
$ws.Hyperlinks.Add($ws.Range("C9"), "mailto:sushantguest1@test.com")
$ws.Hyperlinks.Add($ws.Range("C8"), "mailto:sushantcohost@test.com")

# Hyperlinks.Add re-stamps its own cell style; restore the plain
# "hyperlink-look" style (no border) used by the rest of the e-mail column
$ws.Range("C6").Copy()
$ws.Range("C7").PasteSpecial(-4122)
$ws.Range("C8").PasteSpecial(-4122)
$ws.Range("C9").PasteSpecial(-4122)
$ws.Range("C10").PasteSpecial(-4122)

# --- Leave the selection on the last row that was touched ---
$ws.Range("D10:G10").Select() | Out-Null
